# Update the "2. Data reporter" contact block with the new organization's
# information, and move the active selection to B8 (matching the saved
# workbook's last selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value  = "Kalymbetova Yryskan"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

$ws.Range("B8").Select()
